# Figure 3.xlsx — "Update lite version with new added pictures"
#
# 1) Fix the "ScheduleD DD=1" typo -> "Scheduled DD=1" in the header row
#    (D1) of the three lite-version sheets (Sheet4/Sheet5/Sheet6), which
#    also updates the matching chart series caption.
# 2) Update the view state: Sheet4 becomes the active/selected tab (was
#    Sheet6), with a new zoom level and selected cell; Sheet5 and Sheet6
#    get their own new selected cells while staying inactive.

$wb = $excel.ActiveWorkbook

$fixedText = "Scheduled DD=1"

$sheetNames = @("Sheet4", "Sheet5", "Sheet6")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D1").Value = $fixedText

    # Refresh the chart's cached series caption so the legend/tooltip text
    # matches the corrected header instead of the stale "ScheduleD DD=1".
    $co = $ws.ChartObjects().Item(1)
    $chart = $co.Chart
    $ser = $chart.SeriesCollection().Item(3)
    $ser.Name = $fixedText
}

# Sheet5: new selection, sheet stays in the background.
$ws5 = $wb.Worksheets.Item("Sheet5")
[void]$ws5.Activate()
[void]$ws5.Range("D14").Select()

# Sheet6: new selection, was previously the active tab but no longer is.
$ws6 = $wb.Worksheets.Item("Sheet6")
[void]$ws6.Activate()
[void]$ws6.Range("C1").Select()

# Sheet4: becomes the active tab, gets a new zoom level and selected cell.
$ws4 = $wb.Worksheets.Item("Sheet4")
[void]$ws4.Activate()
[void]$ws4.Range("F10").Select()
$excel.ActiveWindow.Zoom = 115
